$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'31.336.75"
$ws.Range("E2").Value = "  +3.09%  "

$ws.Range("D3").Value = "'2.006.64"
$ws.Range("E3").Value = "  +7.17%  "

$ws.Range("D5").Value = "'0.7779"
$ws.Range("E5").Value = "  +64.96%  "

$ws.Range("D6").Value = "'259.89"
$ws.Range("E6").Value = "  +6.11%  "

$ws.Range("D7").Value = "'1.002"
$ws.Range("E7").Value = "  +0.24%  "

$ws.Range("D8").Value = "'0.3571"
$ws.Range("E8").Value = "  +24.16%  "

$ws.Range("D9").Value = "'28.60"
$ws.Range("E9").Value = "  +30.73%  "

$ws.Range("E10").Value = "  +8.84%  "

$ws.Range("D11").Value = "'0.8584"
$ws.Range("E11").Value = "  +17.64%  "

$ws.Range("D12").Value = "'0.08207"
$ws.Range("E12").Value = "  +5.33%  "

$ws.Range("B13").Value = "Litecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D13").Value = "'101.54"
$ws.Range("E13").Value = "  +1.34%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "'2.004.85"
$ws.Range("E14").Value = "  +7.07%  "

$ws.Range("D15").Value = "'5.634"
$ws.Range("E15").Value = "  +8.87%  "

$ws.Range("D16").Value = "'15.64"
$ws.Range("E16").Value = "  +19.35%  "

$ws.Range("D17").Value = "'274.67"
$ws.Range("E17").Value = "  -4.02%  "

$ws.Range("D18").Value = "'31.335.36"
$ws.Range("E18").Value = "  +3.14%  "

$ws.Range("D19").Value = "'5.948"
$ws.Range("E19").Value = "  +11.41%  "

$ws.Range("D20").Value = "'0.000008006"
$ws.Range("E20").Value = "  +6.79%  "

$ws.Range("D21").Value = "'2.271.84"
$ws.Range("E21").Value = "  +7.52%  "

$ws.Range("D22").Value = "'1.002"
$ws.Range("E22").Value = "  +0.21%  "

$ws.Range("D23").Value = "'1.002"
$ws.Range("E23").Value = "  +0.20%  "

$ws.Range("D24").Value = "'7.175"
$ws.Range("E24").Value = "  +13.50%  "

$ws.Range("D25").Value = "'10.10"
$ws.Range("E25").Value = "  +11.75%  "

$ws.Range("D26").Value = "'166.02"
$ws.Range("E26").Value = "  +1.76%  "

$ws.Range("D27").Value = "'0.1481"
$ws.Range("E27").Value = "  +52.96%  "

$ws.Range("D28").Value = "'20.07"
$ws.Range("E28").Value = "  +5.73%  "

$ws.Range("D29").Value = "'2.395"
$ws.Range("E29").Value = "  +26.22%  "

$ws.Range("E30").Value = "  +9.21%  "

$ws.Range("D31").Value = "'4.638"
$ws.Range("E31").Value = "  +9.64%  "

$ws.Range("D32").Value = "'1.366"
$ws.Range("E32").Value = "  +3.59%  "

$ws.Range("D33").Value = "'4.427"
$ws.Range("E33").Value = "  +6.75%  "

$ws.Range("D34").Value = "'0.05236"
$ws.Range("E34").Value = "  +8.94%  "

$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'1.224"
$ws.Range("E35").Value = "  +8.64%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.7798"
$ws.Range("E36").Value = "  +13.13%  "

$ws.Range("D37").Value = "'2.818"
$ws.Range("E37").Value = "  +3.39%  "

$ws.Range("D38").Value = "'0.02009"
$ws.Range("E38").Value = "  +5.72%  "

$ws.Range("D39").Value = "'2.941"
$ws.Range("E39").Value = "  +3.28%  "

$ws.Range("D40").Value = "'6.730"
$ws.Range("E40").Value = "  +7.12%  "

$ws.Range("D41").Value = "'80.07"
$ws.Range("E41").Value = "  +5.19%  "

$ws.Range("D42").Value = "'0.4749"
$ws.Range("E42").Value = "  +12.39%  "

$ws.Range("D43").Value = "'2.158"
$ws.Range("E43").Value = "  +9.77%  "

$ws.Range("D44").Value = "'107.71"
$ws.Range("E44").Value = "  +6.43%  "

$ws.Range("D45").Value = "'0.8596"
$ws.Range("E45").Value = "  +4.31%  "

$ws.Range("D46").Value = "'1.002"
$ws.Range("E46").Value = "  +0.34%  "

$ws.Range("D47").Value = "'7.833"
$ws.Range("E47").Value = "  +11.58%  "

$ws.Range("D48").Value = "'9.964"
$ws.Range("E48").Value = "  +2.16%  "

$ws.Range("D49").Value = "'0.4366"
$ws.Range("E49").Value = "  +11.63%  "

$ws.Range("D50").Value = "'36.98"
$ws.Range("E50").Value = "  +5.54%  "

$ws.Range("D51").Value = "'0.1203"
$ws.Range("E51").Value = "  +14.56%  "
